$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9..114 down to 10..115
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new record's data
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44545
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 100112045
$ws.Range("G9").Value = "Zapallo"
$ws.Range("H9").Value = "Camote"
$ws.Range("I9").Value = "1a nueva(o)"
$ws.Range("J9").Value = 600
$ws.Range("K9").Value = 600
$ws.Range("L9").Value = 650
$ws.Range("M9").Value = 625
$ws.Range("N9").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O9").Value = "Región de O'Higgins"
$ws.Range("P9").Value = 625
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Hortaliza"
